$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1441441441441441
$ws.Range("C2").Value = 0.6621621621621622
$ws.Range("J2").Value = 0.01351351351351351
$ws.Range("P2").Value = 0.0990990990990991
$ws.Range("S2").Value = 0.08108108108108109

$ws.Range("C3").Value = 0.01333333333333333
$ws.Range("J3").Value = 0.02666666666666667
$ws.Range("P3").Value = 0.7666666666666667
$ws.Range("S3").Value = 0.1933333333333333

$ws.Range("P4").Value = 0.7058823529411765
$ws.Range("S4").Value = 0.2941176470588235

$ws.Range("B6").Value = 0.04736842105263158
$ws.Range("D6").Value = 0.005263157894736842
$ws.Range("F6").Value = 0.03684210526315789
$ws.Range("J6").Value = 0.2421052631578947
$ws.Range("O6").Value = 0.01052631578947368
$ws.Range("Q6").Value = 0.1947368421052632
$ws.Range("R6").Value = 0.1052631578947368
$ws.Range("S6").Value = 0.3578947368421053

$ws.Range("B7").Value = 0.06870229007633588
$ws.Range("D7").Value = 0.007633587786259542
$ws.Range("F7").Value = 0.07633587786259542
$ws.Range("J7").Value = 0.09923664122137404
$ws.Range("O7").Value = 0.03053435114503817
$ws.Range("Q7").Value = 0.1908396946564886
$ws.Range("R7").Value = 0.1297709923664122
$ws.Range("S7").Value = 0.3969465648854962

$ws.Range("B8").Value = 0.07788944723618091
$ws.Range("D8").Value = 0.01507537688442211
$ws.Range("F8").Value = 0.05276381909547739
$ws.Range("J8").Value = 0.135678391959799
$ws.Range("O8").Value = 0.01256281407035176
$ws.Range("Q8").Value = 0.1934673366834171
$ws.Range("R8").Value = 0.1080402010050251
$ws.Range("S8").Value = 0.4045226130653266

$ws.Range("B9").Value = 0.1152073732718894
$ws.Range("D9").Value = 0.009216589861751152
$ws.Range("F9").Value = 0.05069124423963134
$ws.Range("J9").Value = 0.1474654377880184
$ws.Range("O9").Value = 0.009216589861751152
$ws.Range("Q9").Value = 0.1336405529953917
$ws.Range("R9").Value = 0.119815668202765
$ws.Range("S9").Value = 0.4147465437788018

$ws.Range("B10").Value = 0.08682170542635659
$ws.Range("D10").Value = 0.01937984496124031
$ws.Range("F10").Value = 0.06666666666666667
$ws.Range("J10").Value = 0.1372093023255814
$ws.Range("O10").Value = 0.01472868217054264
$ws.Range("Q10").Value = 0.2232558139534884
$ws.Range("R10").Value = 0.1085271317829457
$ws.Range("S10").Value = 0.3434108527131783

$ws.Range("G11").Value = 0.165
$ws.Range("J11").Value = 0.065
$ws.Range("K11").Value = 0.21
$ws.Range("L11").Value = 0.55
$ws.Range("S11").Value = 0.01

$ws.Range("G12").Value = 0.646551724137931
$ws.Range("J12").Value = 0.2413793103448276
$ws.Range("K12").Value = 0.01724137931034483
$ws.Range("L12").Value = 0.0603448275862069
$ws.Range("S12").Value = 0.03448275862068965

$ws.Range("G13").Value = 0.7428571428571429
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.05714285714285714

$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5

$ws.Range("F15").Value = 0.02222222222222222
$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("I15").Value = 0.07222222222222222
$ws.Range("J15").Value = 0.3277777777777778
$ws.Range("K15").Value = 0.06111111111111111
$ws.Range("M15").Value = 0.02222222222222222
$ws.Range("N15").Value = 0.005555555555555556
$ws.Range("O15").Value = 0.08888888888888889
$ws.Range("S15").Value = 0.2333333333333333

$ws.Range("F16").Value = 0.0130718954248366
$ws.Range("H16").Value = 0.1568627450980392
$ws.Range("I16").Value = 0.09803921568627451
$ws.Range("J16").Value = 0.4836601307189543
$ws.Range("K16").Value = 0.08496732026143791
$ws.Range("M16").Value = 0.0196078431372549
$ws.Range("N16").Value = 0.006535947712418301
$ws.Range("O16").Value = 0.06535947712418301
$ws.Range("S16").Value = 0.0718954248366013

$ws.Range("F17").Value = 0.01545253863134658
$ws.Range("H17").Value = 0.2097130242825607
$ws.Range("I17").Value = 0.08167770419426049
$ws.Range("J17").Value = 0.4701986754966888
$ws.Range("K17").Value = 0.05298013245033113
$ws.Range("M17").Value = 0.01324503311258278
$ws.Range("O17").Value = 0.04856512141280353
$ws.Range("S17").Value = 0.108167770419426

$ws.Range("F18").Value = 0.0163265306122449
$ws.Range("H18").Value = 0.1306122448979592
$ws.Range("I18").Value = 0.08571428571428572
$ws.Range("J18").Value = 0.5306122448979592
$ws.Range("K18").Value = 0.04489795918367347
$ws.Range("M18").Value = 0.02040816326530612
$ws.Range("O18").Value = 0.05714285714285714
$ws.Range("S18").Value = 0.1142857142857143

$ws.Range("F19").Value = 0.01232394366197183
$ws.Range("H19").Value = 0.1892605633802817
$ws.Range("I19").Value = 0.1161971830985915
$ws.Range("J19").Value = 0.4049295774647887
$ws.Range("K19").Value = 0.07922535211267606
$ws.Range("M19").Value = 0.01496478873239437
$ws.Range("O19").Value = 0.05721830985915493
$ws.Range("S19").Value = 0.1258802816901408
